$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.249.83"
$ws.Range("E2").Value = "  -0.72%  "
$ws.Range("D3").Value = "2.973.48"
$ws.Range("E3").Value = "  +0.82%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "381.50"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.70%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "102.17"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -3.31%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.543"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -1.09%  "
$ws.Range("E8").Value = "  +0.16%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.589"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -1.21%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.87"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.59%  "
$ws.Range("E11").Value = "  -0.29%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0841"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.14%  "
$ws.Range("D13").Value = "3.444.38"
$ws.Range("E13").Value = "  +1.04%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.16"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.94%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.50"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.31%  "
$ws.Range("D16").Value = "2.972.88"
$ws.Range("E16").Value = "  +1.03%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.996"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +5.29%  "
$ws.Range("D18").Value = "51.206.52"
$ws.Range("E18").Value = "  -0.76%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.25"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -2.55%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.36"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.03%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.79"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -2.66%  "
$ws.Range("D22").Value = "0.0₃0955"
$ws.Range("E22").Value = "  +0.22%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.86"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "260.19"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.92%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.89"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +6.19%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.10"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +12.11%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.47"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +9.05%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.117"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +12.71%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.168"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -2.42%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.12"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.67%  "
$ws.Range("E31").Value = "  -0.05%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "25.86"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.52%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "9.83"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -1.49%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "34.29"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -2.24%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "50.83"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.32%  "
$ws.Range("E36").Value = "  -2.43%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0453"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +5.25%  "
$ws.Range("E38").Value = "  -0.15%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.97"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -2.47%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "16.93"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -2.11%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.56"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.85%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.115"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.78%  "
$ws.Range("E43").Value = "  -2.76%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "122.76"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +2.99%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.42"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -4.35%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.08"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -1.35%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.272"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +3.37%  "
$ws.Range("E48").Value = "  +2.37%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.26"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +1.15%  "
$ws.Range("D50").Value = "2.026.73"
$ws.Range("E50").Value = "  -0.68%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0330"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.69%  "
